# Add a new certification entry (row 65) to the "Control" sheet, mirroring
# the formatting of the previous row (64), and wire up its hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control")

# Duplicate row 64's formatting (styles/borders/number formats) onto row 65
# before writing values, so the new row looks consistent with the table.
$ws.Range("B64:I64").Copy()
$ws.Range("B65:I65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New course entry values.
$ws.Range("B65").Value = "Alura"
$ws.Range("C65").Value = "Databricks: conhecendo a ferramenta"
$ws.Range("D65").Value = 8
$ws.Range("E65").Value = 45435
$ws.Range("F65").Value = "https://cursos.alura.com.br/certificate/f151e8d4-1a19-46ae-b002-ac909dafd7fd"
$ws.Range("G65").Value = "Ok"
$ws.Range("H65").Value = "Ok"
$ws.Range("I65").Value = 45435

# Turn the URL in F65 into a real hyperlink (matching the other rows).
$ws.Hyperlinks.Add($ws.Range("F65"), "https://cursos.alura.com.br/certificate/f151e8d4-1a19-46ae-b002-ac909dafd7fd") | Out-Null

# Re-apply F64's cell style (hyperlink font/border) since adding the
# hyperlink resets to the default "Hyperlink" style/font size.
$ws.Range("F64").Copy()
$ws.Range("F65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to the next empty row, like Excel does after
# data entry.
$ws.Range("B66").Select() | Out-Null
